# Apply timekeeper correction for Doug_Kinsey_2025-12-29.xlsx
# - Client name correction: Smith -> Jones
# - Hours correction: 7.5 -> 8.5 (both the detail row and the SUBTOTAL row)
# - Recalculated dependent totals: Rate 30 -> 65, Total 225 -> 552.5 (detail row)
#   and SUBTOTAL Total 225 -> 552.5
# - SUBTOTAL label text updated to reflect new hours: "Reg: 7.5 / OT: 0" -> "Reg: 8.5 / OT: 0"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - detail row
$ws.Range("B2").Value = "Jones"
$ws.Range("C2").Value = 8.5
$ws.Range("E2").Value = 65
$ws.Range("F2").Value = 552.5

# Row 4 - SUBTOTAL row
$ws.Range("C4").Value = 8.5
$ws.Range("D4").Value = "Reg: 8.5 / OT: 0"
$ws.Range("F4").Value = 552.5

$wb.Save()
